$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Temporarily force the Price (D) / Volume(1h) (E) columns to text format
# so numeric-looking values (e.g. "236.81", "9.240") are written as literal
# text strings rather than being re-interpreted as numbers by Excel -
# matching the original workbook, where these columns hold text values.
$dataRange = $ws.Range("D2:E51")
$dataRange.NumberFormat = "@"

$ws.Range("D2").Value = '30.277.00'
$ws.Range("E2").Value = '  -0.31%  '
$ws.Range("D3").Value = '1.864.45'
$ws.Range("D4").Value = '0.9993'
$ws.Range("E4").Value = '  -0.13%  '
$ws.Range("D5").Value = '236.81'
$ws.Range("E5").Value = '  +0.59%  '
$ws.Range("D6").Value = '0.9993'
$ws.Range("E6").Value = '  -0.10%  '
$ws.Range("D7").Value = '0.4711'
$ws.Range("E7").Value = '  +0.89%  '
$ws.Range("D8").Value = '0.2915'
$ws.Range("E8").Value = '  +2.34%  '
$ws.Range("D9").Value = '0.06559'
$ws.Range("E9").Value = '  -0.20%  '
$ws.Range("D10").Value = '21.96'
$ws.Range("E10").Value = '  +2.64%  '
$ws.Range("D11").Value = '0.07928'
$ws.Range("E11").Value = '  +0.28%  '
$ws.Range("D12").Value = '98.13'
$ws.Range("E12").Value = '  +0.14%  '
$ws.Range("D13").Value = '1.869.07'
$ws.Range("E13").Value = '  -0.13%  '
$ws.Range("D14").Value = '5.168'
$ws.Range("E14").Value = '  +0.74%  '
$ws.Range("E15").Value = '  +0.96%  '
$ws.Range("D16").Value = '267.19'
$ws.Range("E16").Value = '  -4.78%  '
$ws.Range("D17").Value = '30.271.09'
$ws.Range("E17").Value = '  -0.33%  '
$ws.Range("D18").Value = '13.77'
$ws.Range("E18").Value = '  +8.51%  '
$ws.Range("D19").Value = '0.9997'
$ws.Range("E19").Value = '  +0.00%  '
$ws.Range("D20").Value = '0.000007425'
$ws.Range("E20").Value = '  +1.58%  '
$ws.Range("D21").Value = '2.113.61'
$ws.Range("E21").Value = '  -0.05%  '
$ws.Range("D22").Value = '5.314'
$ws.Range("E22").Value = '  -3.52%  '
$ws.Range("E23").Value = '  -0.10%  '
$ws.Range("D24").Value = '6.193'
$ws.Range("E24").Value = '  -0.13%  '
$ws.Range("D25").Value = '167.56'
$ws.Range("E25").Value = '  +1.46%  '
$ws.Range("D26").Value = '9.240'
$ws.Range("E26").Value = '  -0.40%  '
$ws.Range("D27").Value = '18.96'
$ws.Range("E27").Value = '  -1.11%  '
$ws.Range("D28").Value = '1.959'
$ws.Range("E28").Value = '  +1.00%  '
$ws.Range("D29").Value = '1.393'
$ws.Range("E29").Value = '  +1.35%  '
$ws.Range("D30").Value = '0.09870'
$ws.Range("E30").Value = '  +1.40%  '
$ws.Range("D31").Value = '4.385'
$ws.Range("E31").Value = '  -1.02%  '
$ws.Range("D32").Value = '1.472'
$ws.Range("E32").Value = '  -0.39%  '
$ws.Range("D33").Value = '4.058'
$ws.Range("E33").Value = '  -1.59%  '
$ws.Range("D34").Value = '0.04720'
$ws.Range("E34").Value = '  +0.55%  '
$ws.Range("D35").Value = '1.134'
$ws.Range("E35").Value = '  +1.20%  '
$ws.Range("E36").Value = '  -0.47%  '
$ws.Range("D37").Value = '2.704'
$ws.Range("E37").Value = '  -0.36%  '
$ws.Range("E38").Value = '  +1.04%  '
$ws.Range("D39").Value = '2.615'
$ws.Range("E39").Value = '  +2.85%  '
$ws.Range("D40").Value = '6.306'
$ws.Range("E40").Value = '  -0.44%  '
$ws.Range("D41").Value = '74.43'
$ws.Range("E41").Value = '  +1.75%  '
$ws.Range("D42").Value = '1.954'
$ws.Range("E42").Value = '  +0.35%  '
$ws.Range("D43").Value = '0.8456'
$ws.Range("D44").Value = '0.4167'
$ws.Range("E44").Value = '  -0.49%  '
$ws.Range("D45").Value = '0.9987'
$ws.Range("E45").Value = '  -0.17%  '
$ws.Range("D46").Value = '103.46'
$ws.Range("E46").Value = '  -0.52%  '
$ws.Range("D47").Value = '7.188'
$ws.Range("E47").Value = '  -0.33%  '
$ws.Range("D48").Value = '954.40'
$ws.Range("E48").Value = '  +2.17%  '
$ws.Range("D49").Value = '9.239'
$ws.Range("E49").Value = '  +0.93%  '
$ws.Range("D50").Value = '34.16'
$ws.Range("E50").Value = '  +0.07%  '
$ws.Range("D51").Value = '0.05652'
$ws.Range("E51").Value = '  +0.32%  '

# Restore the default cell style so the text format change above does not
# leave a lingering style difference on the cells.
$dataRange.Style = "Normal"
